# Computer crashed and autosaves won't go away yet
# Re-apply the recovered summary-statistics blocks (Kpl @ row22-24, Flow_Lac @ row38-40)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 22 (C22:F22) - cell-line labels
$ws.Range("C22").Value = "HK-2"
$ws.Range("D22").Value = "UMRC6"
$ws.Range("E22").Value = "UOK262"
$ws.Range("F22").Value = "UOK + DIDS"

# Row 23: Kpl averages per cell line (columns B1:B3, B4:B6, B9:B11, B13:B16)
$ws.Range("B23").Value = "Kpl"
$ws.Range("C23").Formula = '=AVERAGE(B$1:B$3)'
$ws.Range("D23").Formula = '=AVERAGE(B$4:B$6)'
$ws.Range("E23").Formula = '=AVERAGE(B$9:B$11)'
$ws.Range("F23").Formula = '=AVERAGE(B$13:B$16)'

# Row 24: Kpl standard error of the mean
$ws.Range("C24").Formula = '=STDEV(B$1:B$3)/SQRT(COUNT(B$1:B$3))'
$ws.Range("D24").Formula = '=STDEV(B$4:B$6)/SQRT(COUNT(B$4:B$6))'
$ws.Range("E24").Formula = '=STDEV(B$9:B$11)/SQRT(COUNT(B$9:B$11))'
$ws.Range("F24").Formula = '=STDEV(B$13:B$16)/SQRT(COUNT(B$13:B$16))'

# Header row 38 (G38:J38) - cell-line labels
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# Row 39: Flow_Lac averages per cell line (columns F1:F3, F4:F6, F9:F11, F13:F16)
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = '=AVERAGE(F$1:F$3)'
$ws.Range("H39").Formula = '=AVERAGE(F$4:F$6)'
$ws.Range("I39").Formula = '=AVERAGE(F$9:F$11)'
$ws.Range("J39").Formula = '=AVERAGE(F$13:F$16)'

# Row 40: Flow_Lac standard error of the mean
$ws.Range("G40").Formula = '=STDEV(F$1:F$3)/SQRT(COUNT(F$1:F$3))'
$ws.Range("H40").Formula = '=STDEV(F$4:F$6)/SQRT(COUNT(F$4:F$6))'
$ws.Range("I40").Formula = '=STDEV(F$9:F$11)/SQRT(COUNT(F$9:F$11))'
$ws.Range("J40").Formula = '=STDEV(F$13:F$16)/SQRT(COUNT(F$13:F$16))'

# Match the recorded selection state (A17:J40, active cell A17)
$ws.Range("A17:J40").Select()
